# "Drop in results from RMI script"
# The EPA moved/renamed their mortality-risk-valuation FAQ page. Update the
# "About" sheet's source-URL cell (B6) and wire up a live hyperlink to the
# new address, replacing the old (now dead) yosemite.epa.gov link.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$newUrl     = "https://www.epa.gov/environmental-economics/mortality-risk-valuation"
$newLoc     = "whatvalue"
$newDisplay = "https://www.epa.gov/environmental-economics/mortality-risk-valuation - whatvalue"
$newText    = "https://www.epa.gov/environmental-economics/mortality-risk-valuation#whatvalue"

$cell = $ws.Range("B6")

# Add the hyperlink (Address + SubAddress split so the relationship target
# is the bare URL and the "#whatvalue" fragment lives in the hyperlink's
# Location attribute), then restore the cell's literal text and its
# original "Hyperlink" cell style afterwards.
$ws.Hyperlinks.Add($cell, $newUrl, $newLoc, "", $newDisplay)
$cell.Value = $newText
$cell.Style = "Hyperlink"
